# Weekly update: insert the latest Haba (Vega Central Mapocho de Santiago) price
# record as a new row at the top of the data (row 49), pushing the existing
# historical rows down by one. This mirrors the "Fruta / hortaliza, semanal"
# commit which prepends the newest observation to each subsetted sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: insert a blank row above the current row 49.
$ws.Rows.Item(49).Insert()

# Populate the new row with the latest reported values.
$ws.Range("A49").Value = 9
$ws.Range("B49").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C49").Value = "Metropolitana"
$ws.Range("D49").Value = 44498
$ws.Range("E49").Value = 13
$ws.Range("F49").Value = 100112026
$ws.Range("G49").Value = "Haba"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 40
$ws.Range("K49").Value = 6000
$ws.Range("L49").Value = 7000
$ws.Range("M49").Value = 6575
$ws.Range("N49").Value = "`$/saco 25 kilos"
$ws.Range("O49").Value = "Región Metropolitana"
$ws.Range("P49").Value = 263
$ws.Range("Q49").Value = 25
$ws.Range("R49").Value = "Hortaliza"
